$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.007179260253906
$ws.Range("B1").Value = 1.481809854507446
$ws.Range("C1").Value = 3.280383348464966
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.454526901245117
